$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section header row (row 25): "Two pointers (Advance)" ---
# Matches the style of the other section-header rows (e.g. row 20 "Two pointers (Basics)")
$ws.Range("B20").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("B25").Value = "Two pointers (Advance)"
$ws.Rows(25).RowHeight = 39

# --- New data row (row 27): "Container With Most Water" ---
# Matches the style/format of the previous data row in that section (row 23)
$ws.Range("A23:E23").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)

$ws.Range("A27").Value = 11
$ws.Range("B27").Value = "Container With Most Water"
$ws.Range("C27").Value = "Medium"
$ws.Range("D27").Value = "Arrays,Two Pointer,Comparison"
$ws.Range("E27").Value = 45695
$ws.Range("E27").NumberFormat = $ws.Range("E23").NumberFormat

$excel.CutCopyMode = 0

# --- Update selection to match the new "active cell" after the edit ---
$ws.Range("D28").Select()
